# Applies the edits described by the commit diff:
#  - date header bumped by one day
#  - division-problem cells replaced with new problems
#
# Each cell's "before" text is unique at the moment we process it,
# so a wdReplaceOne (replace first occurrence only) search over the
# whole document content, executed in document order, is unambiguous
# even though a couple of the "after" strings duplicate "before"
# strings found later in the document.

$d = $word.ActiveDocument

$replacements = @(
    @("2023-10-20 Friday", "2023-10-21 Saturday"),
    @("11÷7=1, 4", "80÷3=26, 2"),
    @("52÷9=5, 7", "76÷2=38, 0"),
    @("52÷7=7, 3", "57÷4=14, 1"),
    @("57÷2=28, 1", "77÷6=12, 5"),
    @("91÷9=10, 1", "67÷7=9, 4"),
    @("19÷3=6, 1", "67÷7=9, 4"),
    @("48÷6=8, 0", "59÷9=6, 5"),
    @("30÷6=5, 0", "82÷8=10, 2"),
    @("27÷7=3, 6", "31÷5=6, 1"),
    @("88÷6=14, 4", "19÷3=6, 1"),
    @("10÷5=2, 0", "64÷9=7, 1"),
    @("22÷3=7, 1", "92÷4=23, 0"),
    @("89÷5=17, 4", "32÷4=8, 0"),
    @("14÷9=1, 5", "96÷3=32, 0"),
    @("36÷9=4, 0", "81÷7=11, 4"),
    @("33÷9=3, 6", "78÷6=13, 0"),
    @("93÷3=31, 0", "68÷6=11, 2"),
    @("28÷6=4, 4", "37÷2=18, 1"),
    @("31÷4=7, 3", "56÷4=14, 0"),
    @("26÷6=4, 2", "94÷4=23, 2"),
    @("83÷7=11, 6", "15÷7=2, 1"),
    @("39÷8=4, 7", "81÷7=11, 4"),
    @("38÷8=4, 6", "77÷3=25, 2"),
    @("71÷9=7, 8", "78÷7=11, 1"),
    @("42÷2=21, 0", "16÷8=2, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 1)
}
